$d = $word.ActiveDocument

# --- Date heading ---
$d.Content.Find.Execute("2024-01-21 Sunday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-01-22 Monday", 2)

# --- Table cell replacements ---
# Most old values are unique across the document, so a straightforward
# Find/Replace works for them. Two cells share the same original text
# ("95÷6=15, 5") but map to different new values, so those two are
# addressed directly via the Tables/Cell collection instead.

$t = $d.Tables.Item(1)
$t.Cell(9, 2).Range.Text  = "37÷5=7, 2"
$t.Cell(13, 2).Range.Text = "24÷6=4, 0"

$pairs = @(
    @("15÷8=1, 7",  "43÷8=5, 3"),
    @("77÷9=8, 5",  "64÷8=8, 0"),
    @("97÷8=12, 1", "89÷8=11, 1"),
    @("73÷9=8, 1",  "72÷5=14, 2"),
    @("30÷5=6, 0",  "54÷8=6, 6"),
    @("86÷8=10, 6", "26÷3=8, 2"),
    @("44÷8=5, 4",  "65÷2=32, 1"),
    @("96÷2=48, 0", "93÷6=15, 3"),
    @("20÷3=6, 2",  "89÷4=22, 1"),
    @("54÷6=9, 0",  "25÷9=2, 7"),
    @("66÷7=9, 3",  "75÷2=37, 1"),
    @("64÷6=10, 4", "37÷8=4, 5"),
    @("49÷5=9, 4",  "67÷2=33, 1"),
    @("13÷7=1, 6",  "15÷5=3, 0"),
    @("10÷7=1, 3",  "85÷4=21, 1"),
    @("87÷9=9, 6",  "24÷9=2, 6"),
    @("13÷6=2, 1",  "26÷6=4, 2"),
    @("87÷4=21, 3", "28÷4=7, 0"),
    @("59÷4=14, 3", "89÷4=22, 1"),
    @("32÷5=6, 2",  "85÷8=10, 5"),
    @("98÷6=16, 2", "75÷5=15, 0"),
    @("38÷3=12, 2", "63÷7=9, 0"),
    @("34÷4=8, 2",  "80÷2=40, 0")
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair[0], $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair[1], 2)
}
